$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)

# Leetcode number for the existing row 2 ("罗马数字转整数" problem) was left
# blank before; the commit fills it in.
$ws.Range("B2").Value2 = 168

# Insert a new data row (row 3) for the new "Excel column title" / "binary
# sum" problems, inheriting the formatting (style, wrap, etc.) of row 2.
$ws.Rows.Item(3).Insert()

$ws.Range("A3").Value2 = 2
$ws.Range("B3").Value2 = "171`n#1290"
$ws.Range("D3").Value2 = "N进制转换成十进制之和：`n         1 从左到右遍历每一位，即从高位累加，而非低位`n         2 累加和sum，sum = sum * N + val`n         3 原理介绍：`n              例如数字abc，十进制是a*100 + b * 10 + c`n              转换后变成： (a*10+b)*10+c`n              即可循环迭代每一位数字`n"
$ws.Range("C3").Value2 = "#171`n给定一个Excel表格中的列名称，返回其相应的列序号。 `n 例如， `n     A -> 1`n    B -> 2`n    C -> 3`n    ...`n    Z -> 26`n    AA -> 27`n    AB -> 28 `n    ...`n 输入: ""A""  输出: 1`n 输入: ""AB""   输出: 28`n 输入: ""ZY""  输出: 701 `n#1290`n给你一个单链表的引用结点 head。链表中每个结点的值不是 0 就 1。已知此链表是一个整数数字的二进制表示形式。 `n 请你返回该链表所表示数字的 十进制值 。 `n 输入：head = [1,0,1]`n输出：5`n解释：二进制数 (101) 转化为十进制数 (5)"
$ws.Range("E3").Value2 = "数学`n进制转换`n累加"
$ws.Range("F3").Value2 = $ws.Range("F2").Value2
$ws.Range("G3").Value2 = $ws.Range("G2").Value2

$ws.Rows.Item(3).RowHeight = 409.6
